# TestcaseSignUp.xlsx — "đăng kí (hoàn thành)"
# Replace the row-6 test data (a stale "user202 / adfghjkfghjkfgh1234567 ..." row)
# with the final long-input test case values, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A6").Value = "qưertyuiopasdfghjklzxcvbnm"
$ws.Range("B6").Value = "1234567890123456"
$ws.Range("C6").Value = "Kiều Oanh qưertyuiopasdfghjklzxc"
$ws.Range("D6").Value = "qwertyuiopasdfghjklzxcvbnmqwertyuiopasdfghjklzxcvbnm@gmail.com"

$ws.Range("E12").Select()
